$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# 1. Overview sheet: handback status text changes for both locales
#    ("In Translation" -> "Handed back: in sync with en-US")
# ---------------------------------------------------------------------------
$overview = $wb.Worksheets.Item("Overview")
$overview.Range("E2").Value = "Handed back: in sync with en-US"
$overview.Range("F2").Value = "Handed back: in sync with en-US"
$overview.Range("E3").Value = "Handed back: in sync with en-US"
$overview.Range("F3").Value = "Handed back: in sync with en-US"

# Widen the status columns now that the text is longer.
$overview.Columns.Item(5).ColumnWidth = 29.9777047293527
$overview.Columns.Item(6).ColumnWidth = 29.9777047293527

# The "Status" column on each locale sheet shares the very same string
# ("In Translation"); update it there too so every occurrence reflects the
# handback.
$wb.Worksheets.Item("zh-cn").Range("C2").Value = "Handed back: in sync with en-US"
$wb.Worksheets.Item("zh-cn").Range("C3").Value = "Handed back: in sync with en-US"
$wb.Worksheets.Item("de-de").Range("C2").Value = "Handed back: in sync with en-US"
$wb.Worksheets.Item("de-de").Range("C3").Value = "Handed back: in sync with en-US"

# ---------------------------------------------------------------------------
# 2. zh-cn sheet: fill in target/handback file + handback datetime for the
#    two rows, and add hyperlinks that mirror the "Source File Name" links.
# ---------------------------------------------------------------------------
$zhcn = $wb.Worksheets.Item("zh-cn")

$zhcn.Range("I2").Value = "6eea64ae-6632-43f2-969c-4fd297c51373.md"
$zhcn.Range("J2").Value = "6eea64ae-6632-43f2-969c-4fd297c51373.dfb27053638d424c5df6d9311dec3553299e6493.zh-cn.xlf"
$zhcn.Range("K2").Value = "2016-08-12 02:47:23"

$zhcn.Range("I3").Value = "e84e0925-2728-4339-93c6-04e674dbfa31.md"
$zhcn.Range("J3").Value = "e84e0925-2728-4339-93c6-04e674dbfa31.ecb508a810c614b8b8494d00d432b91f64f5e959.zh-cn.xlf"
$zhcn.Range("K3").Value = "2016-08-12 02:47:23"

# Re-create the A3 hyperlink after the new I2 hyperlink so relationship ids
# land in the same order as row-by-row link generation (A2, I2, A3, I3).
foreach ($hl in $zhcn.Hyperlinks) {
  $addr = $hl.Range.Address()
  if ($addr -eq '$A$3') {
    $hl.Delete()
  }
}

$zhcn.Hyperlinks.Add($zhcn.Range("I2"), "https://github.com/OpenLocalizationTestOrg/oltest/blob/8c602ded64e536251a9a21c99bec246dcb3411cc/e2e/6eea64ae-6632-43f2-969c-4fd297c51373.md", "", "", "6eea64ae-6632-43f2-969c-4fd297c51373.md")
$zhcn.Hyperlinks.Add($zhcn.Range("A3"), "https://github.com/OpenLocalizationTestOrg/oltest/blob/8c602ded64e536251a9a21c99bec246dcb3411cc/e2e/e84e0925-2728-4339-93c6-04e674dbfa31.md", "", "", "e84e0925-2728-4339-93c6-04e674dbfa31.md")
$zhcn.Hyperlinks.Add($zhcn.Range("I3"), "https://github.com/OpenLocalizationTestOrg/oltest/blob/8c602ded64e536251a9a21c99bec246dcb3411cc/e2e/e84e0925-2728-4339-93c6-04e674dbfa31.md", "", "", "e84e0925-2728-4339-93c6-04e674dbfa31.md")

$zhcn.Columns.Item(3).ColumnWidth = 29.9777047293527
$zhcn.Columns.Item(9).ColumnWidth = 40
$zhcn.Columns.Item(10).ColumnWidth = 40

# ---------------------------------------------------------------------------
# 3. de-de sheet: same shape of edit, different handback datetime value.
# ---------------------------------------------------------------------------
$dede = $wb.Worksheets.Item("de-de")

$dede.Range("I2").Value = "6eea64ae-6632-43f2-969c-4fd297c51373.md"
$dede.Range("J2").Value = "6eea64ae-6632-43f2-969c-4fd297c51373.dfb27053638d424c5df6d9311dec3553299e6493.de-de.xlf"
$dede.Range("K2").Value = "2016-08-12 02:47:31"

$dede.Range("I3").Value = "e84e0925-2728-4339-93c6-04e674dbfa31.md"
$dede.Range("J3").Value = "e84e0925-2728-4339-93c6-04e674dbfa31.ecb508a810c614b8b8494d00d432b91f64f5e959.de-de.xlf"
$dede.Range("K3").Value = "2016-08-12 02:47:31"

foreach ($hl in $dede.Hyperlinks) {
  $addr = $hl.Range.Address()
  if ($addr -eq '$A$3') {
    $hl.Delete()
  }
}

$dede.Hyperlinks.Add($dede.Range("I2"), "https://github.com/OpenLocalizationTestOrg/oltest/blob/8c602ded64e536251a9a21c99bec246dcb3411cc/e2e/6eea64ae-6632-43f2-969c-4fd297c51373.md", "", "", "6eea64ae-6632-43f2-969c-4fd297c51373.md")
$dede.Hyperlinks.Add($dede.Range("A3"), "https://github.com/OpenLocalizationTestOrg/oltest/blob/8c602ded64e536251a9a21c99bec246dcb3411cc/e2e/e84e0925-2728-4339-93c6-04e674dbfa31.md", "", "", "e84e0925-2728-4339-93c6-04e674dbfa31.md")
$dede.Hyperlinks.Add($dede.Range("I3"), "https://github.com/OpenLocalizationTestOrg/oltest/blob/8c602ded64e536251a9a21c99bec246dcb3411cc/e2e/e84e0925-2728-4339-93c6-04e674dbfa31.md", "", "", "e84e0925-2728-4339-93c6-04e674dbfa31.md")

$dede.Columns.Item(3).ColumnWidth = 29.9777047293527
$dede.Columns.Item(9).ColumnWidth = 40
$dede.Columns.Item(10).ColumnWidth = 40

Write-Host "Handback report generated"
